# Commit: update the up and down data
#
# This script reproduces, via the Excel COM object model, the edit that:
#  - populates new "result" columns (C/D/E) on the "down" and "up" sheets
#  - adds the corresponding sheet-local defined names (res/res_1/res_2)
#  - makes the "down" sheet the active tab/selection
#  - sizes the new columns

$wb = $excel.ActiveWorkbook

$wsDown = $wb.Worksheets.Item("down")
$wsUp   = $wb.Worksheets.Item("up")

# ---------------------------------------------------------------------------
# 1. Populate "down" sheet columns C (QM+phn), D (X6), E (VopX)
# ---------------------------------------------------------------------------
$wsDown.Range("C1").Value = "QM+phn"
$wsDown.Range("D1").Value = "X6"
$wsDown.Range("E1").Value = "VopX"
$wsDown.Range("C2").Value = -1.663976472
$wsDown.Range("D2").Value = -1.5973
$wsDown.Range("E2").Value = -1.5455
$wsDown.Range("C3").Value = -1.674777878
$wsDown.Range("D3").Value = -1.6272
$wsDown.Range("E3").Value = -1.5806
$wsDown.Range("C4").Value = -1.686285249
$wsDown.Range("D4").Value = -1.6537
$wsDown.Range("E4").Value = -1.6139
$wsDown.Range("C5").Value = -1.695119714
$wsDown.Range("D5").Value = -1.6764
$wsDown.Range("E5").Value = -1.6448
$wsDown.Range("C6").Value = -1.701789497
$wsDown.Range("D6").Value = -1.6944
$wsDown.Range("E6").Value = -1.6729
$wsDown.Range("C7").Value = -1.705037084
$wsDown.Range("D7").Value = -1.707
$wsDown.Range("E7").Value = -1.6975
$wsDown.Range("C8").Value = -1.704158794
$wsDown.Range("D8").Value = -1.7133
$wsDown.Range("E8").Value = -1.718
$wsDown.Range("C9").Value = -1.699050211
$wsDown.Range("D9").Value = -1.7125
$wsDown.Range("E9").Value = -1.7337
$wsDown.Range("C10").Value = -1.688868799
$wsDown.Range("D10").Value = -1.7033
$wsDown.Range("E10").Value = -1.7437
$wsDown.Range("C11").Value = -1.672887533
$wsDown.Range("D11").Value = -1.6847
$wsDown.Range("E11").Value = -1.747
$wsDown.Range("C12").Value = -1.650311357
$wsDown.Range("D12").Value = -1.6553
$wsDown.Range("E12").Value = -1.7424
$wsDown.Range("C13").Value = -1.620320181
$wsDown.Range("D13").Value = -1.6136
$wsDown.Range("E13").Value = -1.7288
$wsDown.Range("C14").Value = -1.582144172
$wsDown.Range("D14").Value = -1.5582
$wsDown.Range("E14").Value = -1.7047
$wsDown.Range("C15").Value = -1.534574247
$wsDown.Range("D15").Value = -1.4871
$wsDown.Range("E15").Value = -1.6684
$wsDown.Range("C16").Value = -1.476615551
$wsDown.Range("D16").Value = -1.3984
$wsDown.Range("E16").Value = -1.6181

# ---------------------------------------------------------------------------
# 2. Populate "up" sheet columns C (QM+phn), D (X6), E (VopX)
# ---------------------------------------------------------------------------
$wsUp.Range("C1").Value = "QM+phn"
$wsUp.Range("D1").Value = "X6"
$wsUp.Range("E1").Value = "VopX"
$wsUp.Range("C2").Value = -0.663976472
$wsUp.Range("D2").Value = -0.6423
$wsUp.Range("E2").Value = -0.7287
$wsUp.Range("C3").Value = -0.674777878
$wsUp.Range("D3").Value = -0.6592
$wsUp.Range("E3").Value = -0.7327
$wsUp.Range("C4").Value = -0.686285249
$wsUp.Range("D4").Value = -0.6745
$wsUp.Range("E4").Value = -0.7337
$wsUp.Range("C5").Value = -0.695119714
$wsUp.Range("D5").Value = -0.6877
$wsUp.Range("E5").Value = -0.7312
$wsUp.Range("C6").Value = -0.701789497
$wsUp.Range("D6").Value = -0.6982
$wsUp.Range("E6").Value = -0.7247
$wsUp.Range("C7").Value = -0.705037084
$wsUp.Range("D7").Value = -0.7055
$wsUp.Range("E7").Value = -0.7136
$wsUp.Range("C8").Value = -0.704158794
$wsUp.Range("D8").Value = -0.7087
$wsUp.Range("E8").Value = -0.6971
$wsUp.Range("C9").Value = -0.699050211
$wsUp.Range("D9").Value = -0.707
$wsUp.Range("E9").Value = -0.6745
$wsUp.Range("C10").Value = -0.688868799
$wsUp.Range("D10").Value = -0.6992
$wsUp.Range("E10").Value = -0.6449
$wsUp.Range("C11").Value = -0.672887533
$wsUp.Range("D11").Value = -0.6841
$wsUp.Range("E11").Value = -0.6072
$wsUp.Range("C12").Value = -0.650311357
$wsUp.Range("D12").Value = -0.6603
$wsUp.Range("E12").Value = -0.5604
$wsUp.Range("C13").Value = -0.620320181
$wsUp.Range("D13").Value = -0.6261
$wsUp.Range("E13").Value = -0.5031
$wsUp.Range("C14").Value = -0.582144172
$wsUp.Range("D14").Value = -0.5796
$wsUp.Range("E14").Value = -0.4339
$wsUp.Range("C15").Value = -0.534574247
$wsUp.Range("D15").Value = -0.5183
$wsUp.Range("E15").Value = -0.3512
$wsUp.Range("C16").Value = -0.476615551
$wsUp.Range("D16").Value = -0.4397
$wsUp.Range("E16").Value = -0.2531

# ---------------------------------------------------------------------------
# 3. Column widths for the new columns
# ---------------------------------------------------------------------------
$wsDown.Columns.Item(3).ColumnWidth = 11.830729166666666
$wsDown.Columns.Item(4).ColumnWidth = 6.830729166666667
$wsDown.Columns.Item(5).ColumnWidth = 6.830729166666667

$wsUp.Columns.Item(3).ColumnWidth = 11.830729166666666
$wsUp.Columns.Item(4).ColumnWidth = 6.830729166666667
$wsUp.Columns.Item(5).ColumnWidth = 6.830729166666667

# ---------------------------------------------------------------------------
# 4. Sheet-local defined names: res / res_1 / res_2 for "down" and "up"
# ---------------------------------------------------------------------------
$wsDown.Names.Add("res",   "=down!`$C`$2:`$C`$16")
$wsDown.Names.Add("res_1", "=down!`$D`$2:`$D`$16")
$wsDown.Names.Add("res_2", "=down!`$E`$2:`$E`$16")

$wsUp.Names.Add("res",   "=up!`$C`$2:`$C`$16")
$wsUp.Names.Add("res_1", "=up!`$E`$2:`$E`$16")
$wsUp.Names.Add("res_2", "=up!`$D`$2:`$D`$16")

# ---------------------------------------------------------------------------
# 5. Selection / active tab: "down" becomes the active sheet & tab
# ---------------------------------------------------------------------------
$wsUp.Activate()
$wsUp.Range("A1:E16").Select()

$wsDown.Activate()
$wsDown.Range("A1:E16").Select()
